$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 -> becomes old Row 4's data
$ws.Range("D3").Value = 44874
$ws.Range("M3").Value = 67
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 16000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("R3").Value = "Provincia de Los Andes"
$ws.Range("S3").Value = 1600
$ws.Range("T3").Value = 10

# Row 4 -> becomes old Row 5's data
$ws.Range("D4").Value = 45222
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 1500
$ws.Range("T4").Value = 10

# Row 5 -> becomes old Row 3's data
$ws.Range("D5").Value = 44855
$ws.Range("M5").Value = 25
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = "$/bandeja 5 kilos"
$ws.Range("R5").Value = "Provincia de Los Andes"
$ws.Range("S5").Value = 3000
$ws.Range("T5").Value = 5
